# Parse geocoding responses into real street addresses (was placeholder
# letters A-F) and tidy up the sheet so it's ready for the request map.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "100 park ave west, mansfield OH"
$ws.Range("A3").Value = "200 park ave west, mansfield OH"
$ws.Range("A4").Value = "580 woodland road, mansfield OH"
$ws.Range("A6").Value = "397 west 4th, mansfield OH"
$ws.Range("A7").Value = "3366 muskie drive, mansfield OH"
$ws.Range("A5").Value = "601 woodland road, mansfield OH"

$ws.Columns.Item(1).ColumnWidth = 19.2

$ws.Range("A5").Select()
